$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 17 updates
$ws.Range("S17").Value = -0.1
$ws.Range("T17").Value = -0.1
$ws.Range("U17").Value = 0
$ws.Range("V17").Value = 0.2
$ws.Range("W17").Value = -0.2

# Row 18 updates
$ws.Range("T18").Value = 0
$ws.Range("U18").Value = 0

# Row 23 updates
$ws.Range("S23").Value = 0.1
$ws.Range("T23").Value = 0.1
$ws.Range("U23").Value = 0.1
$ws.Range("V23").Value = 0.2
$ws.Range("W23").Value = 0.3
$ws.Range("X23").Value = 0.4
$ws.Range("Y23").Value = 1

# Restore the active cell selection to match the saved workbook state
$ws.Range("X17").Select()
